# "adding cms spread cap floor"
#
# The InflationInput.xlsx sheet "Tabelle1" held a small lookup table
# (rows 2-8, column B) describing one particular inflation index (hicp /
# EUR Inflation 31122019 / 105.02 / EUR Real Vol / EUR Inflation Vol /
# 0.03) together with a reviewer's comment on the starting index cell.
# That whole column B of sample/demo values (and the comment that went
# with it) is cleared out, leaving only the label column (A) and the
# empty, but still percentage-formatted, B8 cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Remove the reviewer comment that was anchored on B5 ("...index at 1th
# sept 19...") before clearing the cell it sits on.
if ($ws.Range("B5").Comment -ne $null) {
    $ws.Range("B5").Comment.Delete() | Out-Null
}

# Wipe the sample values out of column B for the descriptive rows (Name,
# Currency, ZC Inflation, Starting Inflation Index, Volatility Real IR,
# Volatility Inflation Index) - the cells disappear entirely, same as
# pressing Delete on a plain unstyled cell.
$ws.Range("B2:B7").ClearContents() | Out-Null

# B8 (Mean reversion Real IR) keeps its percentage number format (style
# index 1) but the 0.03 value itself is cleared, leaving an empty but
# still-formatted cell.
$ws.Range("B8").ClearContents() | Out-Null

# Reflect the new selection left behind on the sheet.
$ws.Activate() | Out-Null
$ws.Range("B17").Select() | Out-Null
